$d = $word.ActiveDocument

$find = "Dates de la campanya Constel·lació de Taure 2022: 16-25 de gener"
$replace = "Dates de la campanya 2022 en què usem la constel·lació, Constel·lació de Taure 16-25 de gener"

$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
